# Applies the "Added many more features" revision to the Miss Midas review.
#
# Summary of changes:
#  1. Title (Heading1 + the bold repeat near the end) changed.
#  2. "What we like" bullet list: reworded first bullet, inserted two new
#     bullets after it, and removed two bullets at the end (net: still four
#     bullets, but different content/order).
#  3. "What we don't like" bullet list: the two bullets effectively swapped
#     order (each became the other's text).
#  4. Meta-description (italic) paragraph reworded.

$d = $word.ActiveDocument

function New-PPkgXml([string]$text, [string]$pStyle, [bool]$listFormatting) {
  # Builds a FlatOPC single-paragraph package that InsertXML understands,
  # reproducing the paragraph shape used throughout this document: an empty
  # run followed by the text run.
  $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
  $pPr = '<w:pStyle w:val="' + $pStyle + '"/>'
  if ($listFormatting) {
    $pPr = $pPr + '<w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/>'
  }
  return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:pPr>' + $pPr + '</w:pPr><w:r/><w:r><w:t>' + $escaped + '</w:t></w:r></w:p></w:body>' +
    '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1. Title change (replaces both occurrences: the Heading1 at the top and
#    the bold run repeated further down).
# ---------------------------------------------------------------------
[void]$d.Content.Find.Execute(
  "Play Miss Midas Free - Review of Miss Midas Slot Game", $true, $false, $false,
  $false, $false, $true, 1, $false,
  "Play Miss Midas for Free - Review and Gameplay Overview", 2)

# ---------------------------------------------------------------------
# 2. "What we like" list.
#    Before: Golden Touch bonus feature / Unique female protagonist /
#            Beautifully designed symbols / Exciting free spins
#    After:  Simple gameplay mechanics / Golden Touch feature /
#            Visually stunning design / Unique female protagonist
# ---------------------------------------------------------------------

# Find the "Golden Touch bonus feature" bullet paragraph.
$likeFirst = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Golden Touch bonus feature") {
    $likeFirst = $d.Paragraphs.Item($i)
    break
  }
}

# Reword it in place (keeps the empty-run + text-run paragraph shape).
[void]$likeFirst.Range.InsertXML((New-PPkgXml "Simple gameplay mechanics" "ListBullet" $true))

# Insert the two new bullets right after it, in order.
[void]$likeFirst.Range.InsertXML((New-PPkgXml "Visually stunning design" "ListBullet" $true), "After")
[void]$likeFirst.Range.InsertXML((New-PPkgXml "Golden Touch feature" "ListBullet" $true), "After")

# Remove the two now-trailing bullets that are no longer present
# ("Beautifully designed symbols" and "Exciting free spins"), which sit
# right after "Unique female protagonist".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
  if ($t -eq "Beautifully designed symbols") {
    [void]$d.Paragraphs.Item($i).Range.Delete()
    break
  }
}
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
  if ($t -eq "Exciting free spins") {
    [void]$d.Paragraphs.Item($i).Range.Delete()
    break
  }
}

# ---------------------------------------------------------------------
# 3. "What we don't like" list: swap the two bullets' text.
#    Before: Relatively low RTP / No multiple bonus features
#    After:  No multiple bonus games / Relatively low RTP
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
  if ($t -eq "Relatively low RTP") {
    [void]$d.Paragraphs.Item($i).Range.InsertXML((New-PPkgXml "No multiple bonus games" "ListBullet" $true))
    break
  }
}
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
  if ($t -eq "No multiple bonus features") {
    [void]$d.Paragraphs.Item($i).Range.InsertXML((New-PPkgXml "Relatively low RTP" "ListBullet" $true))
    break
  }
}

# ---------------------------------------------------------------------
# 4. Meta description (italic) paragraph reworded.
# ---------------------------------------------------------------------
[void]$d.Content.Find.Execute(
  "Play Miss Midas slot game for free and read our comprehensive review of its gameplay mechanics, theme, symbol design, and jackpot & RTP. Get ready for a unique adventure!",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Read our review of Miss Midas, a visually stunning slot game with the Golden Touch feature. Play for free and experience the adventure!", 2)
